$wb = $excel.ActiveWorkbook

# The first worksheet ("testing competition #1") is renamed to "Name".
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Name"

# The title cell (B1) holds the competition name and mirrors the sheet name.
$ws.Range("B1").Value = "Name"

# Column B was resized (closest reachable width to 13.24609375 given
# the engine's internal 1/6-character rounding of ColumnWidth).
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334

# The competition date (B3) moved from 12/3/2021 to 12/21/2021.
$ws.Range("B3").Value = 44551.0
